$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values
$ws.Range("B2").Value = 60.780889897022746
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 55.357330951788185
$ws.Range("E2").Value = 68.559182613534276

# Row 3 values
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 56.628916748618288
$ws.Range("D3").Value = 56.621218705819487
$ws.Range("E3").Value = 83.562070323499313

$ws.Range("B1:E3").Select()
